# feat: add 2022-Q1 data
#
# - Inserts a new "2022-Q1" worksheet (fund holdings detail) positioned
#   between the existing "2021-Q1" sheet and the "总计" (totals) sheet.
# - Adds a new summary row for "2022-Q1" at the top of the "总计" sheet's
#   data (pushing the existing "2021-Q1" summary row down by one).

$wb = $excel.ActiveWorkbook

# Helper: write a value into a cell while forcing it to stay TEXT even when
# it looks numeric (e.g. "011351" or "0.0070"), the way the rest of this
# workbook stores its "numeric-looking" descriptive columns. Toggling the
# number format to Text and back to General keeps the cell's stored type as
# text while leaving no lingering cell style behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Helper: apply the bold / centered / thin-bordered look used for the
# row-index column (A) and header rows elsewhere in this workbook.
function Set-IndexStyle($range) {
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108   # xlCenter
    $range.VerticalAlignment = -4160     # xlTop
    $range.Borders.LineStyle = 1         # xlContinuous
}

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q1" sheet right before "总计" so the final tab
#    order is: 2021-Q1, 2022-Q1, 总计
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Header row (row 1) — bold, centered, thin-bordered, like the other sheets.
Set-IndexStyle $q1.Range("B1:H1")

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Row 2 — first fund holding.
Set-IndexStyle $q1.Range("A2")
$q1.Range("A2").Value = 0

Set-TextValue $q1.Range("B2") "011351"
Set-TextValue $q1.Range("C2") "金鹰年年邮益一年持有期混合A"
Set-TextValue $q1.Range("D2") "9.03"
Set-TextValue $q1.Range("E2") "37.02"
Set-TextValue $q1.Range("F2") "1.19"
Set-TextValue $q1.Range("G2") "0.1075"
$q1.Range("H2").Value = 5

# Row 3 — second fund holding (share class C of the same fund).
Set-IndexStyle $q1.Range("A3")
$q1.Range("A3").Value = 1

Set-TextValue $q1.Range("B3") "011352"
Set-TextValue $q1.Range("C3") "金鹰年年邮益一年持有期混合C"
Set-TextValue $q1.Range("D3") "0.59"
Set-TextValue $q1.Range("E3") "37.02"
Set-TextValue $q1.Range("F3") "1.19"
Set-TextValue $q1.Range("G3") "0.0070"
$q1.Range("H3").Value = 5

# ---------------------------------------------------------------------------
# 2. Add the 2022-Q1 summary row to the "总计" sheet, just above the
#    existing 2021-Q1 summary row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# The insert copies row 1's (header) formatting onto the new row; strip it
# back down to plain cells first, then re-apply only what column A needs.
$total.Range("A2:D2").ClearFormats()

Set-IndexStyle $total.Range("A2")
$total.Range("A2").Value = 0

$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.11

# The index column (A) is a sequential 0-based row index (pandas-style),
# so the pre-existing "2021-Q1" row — now pushed down to row 3 — needs its
# index renumbered from 0 to 1.
$total.Range("A3").Value = 1

# ---------------------------------------------------------------------------
# 3. Restore the originally active sheet/selection (unchanged by this edit).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q1").Activate()
